$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J116").Value = 5195.5
$ws.Range("L116").Value = 5195.5
$ws.Range("N116").Value = -12079.5
$ws.Range("H132").Value = 4348.7896
$ws.Range("I132").Value = 1185.6897
$ws.Range("K132").Value = 3557.0691
$ws.Range("M132").Value = -1027.0691
$ws.Range("H137").Value = 10101516
$ws.Range("I137").Value = 558882.8
$ws.Range("J137").Value = 20836978
$ws.Range("K137").Value = 1676648.4
$ws.Range("L137").Value = 62510934
$ws.Range("M137").Value = -1674098.4
$ws.Range("N137").Value = -62516034

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 998.2727
$ws.Range("I2").Value = 931.5238000000001
$ws.Range("K2").Value = 931.5238000000001
$ws.Range("M2").Value = -818.5238000000001
$ws.Range("H5").Value = 908.1667
$ws.Range("I5").Value = 908.1667
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 908.1667
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -796.1667
$ws.Range("N5").ClearContents()
$ws.Range("H45").Value = 105181.45
$ws.Range("I45").Value = 142755.88
$ws.Range("J45").Value = 4983
$ws.Range("K45").Value = 142755.88
$ws.Range("L45").Value = 4983
$ws.Range("M45").Value = -142378.88
$ws.Range("N45").Value = -5737
$ws.Range("H61").Value = 3850.8
$ws.Range("I61").Value = 2405.3914
$ws.Range("K61").Value = 2405.3914
$ws.Range("M61").Value = -2193.3914
$ws.Range("H74").Value = 31251420
$ws.Range("I74").Value = 41667560
$ws.Range("K74").Value = 41667560
$ws.Range("M74").Value = -41666686
$ws.Range("H77").Value = 31251420
$ws.Range("I77").Value = 41667560
$ws.Range("K77").Value = 208337800
$ws.Range("M77").Value = -208333432
$ws.Range("H116").Value = 998.2727
$ws.Range("I116").Value = 931.5238000000001
$ws.Range("K116").Value = 931.5238000000001
$ws.Range("M116").Value = 1362.4762
$ws.Range("H136").Value = 3850.8
$ws.Range("I136").Value = 2405.3914
$ws.Range("K136").Value = 7216.174199999999
$ws.Range("M136").Value = -4666.174199999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 998.2727
$ws.Range("I3").Value = 931.5238000000001
$ws.Range("K3").Value = 931.5238000000001
$ws.Range("M3").Value = -817.5238000000001
$ws.Range("H4").Value = 908.1667
$ws.Range("I4").Value = 908.1667
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 908.1667
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -793.1667
$ws.Range("N4").ClearContents()
$ws.Range("H107").Value = 2351.4211
$ws.Range("I107").Value = 2134.1924
$ws.Range("J107").Value = 2822.0833
$ws.Range("K107").Value = 2134.1924
$ws.Range("L107").Value = 2822.0833
$ws.Range("M107").Value = -214.1923999999999
$ws.Range("N107").Value = -6662.0833
$ws.Range("H134").Value = 2506.8572
$ws.Range("I134").Value = 2134.8276
$ws.Range("J134").Value = 3046.3
$ws.Range("K134").Value = 6404.4828
$ws.Range("L134").Value = 9138.900000000001
$ws.Range("M134").Value = -3869.4828
$ws.Range("N134").Value = -14208.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19234028
$ws.Range("I31").Value = 20836266
$ws.Range("K31").Value = 20836266
$ws.Range("M31").Value = -20835971
$ws.Range("H34").Value = 19234028
$ws.Range("I34").Value = 20836266
$ws.Range("K34").Value = 20836266
$ws.Range("M34").Value = -20836064
$ws.Range("H74").Value = 54999.5
$ws.Range("J74").Value = 54999.5
$ws.Range("L74").Value = 54999.5
$ws.Range("N74").Value = -56747.5
$ws.Range("H77").Value = 54999.5
$ws.Range("J77").Value = 54999.5
$ws.Range("L77").Value = 164998.5
$ws.Range("N77").Value = -173734.5
$ws.Range("H107").Value = 827.3333
$ws.Range("I107").Value = 575.125
$ws.Range("J107").Value = 919.0454999999999
$ws.Range("K107").Value = 575.125
$ws.Range("L107").Value = 919.0454999999999
$ws.Range("M107").Value = 1344.875
$ws.Range("N107").Value = -4759.0455
$ws.Range("H132").Value = 33334976
$ws.Range("I132").Value = 43012132
$ws.Range("J132").Value = 2556.3333
$ws.Range("K132").Value = 129036396
$ws.Range("L132").Value = 7668.999899999999
$ws.Range("M132").Value = -129033866
$ws.Range("N132").Value = -12728.9999
$ws.Range("H134").Value = 2157.7407
$ws.Range("I134").Value = 1930.3334
$ws.Range("K134").Value = 5791.0002
$ws.Range("M134").Value = -3256.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2246.7646
$ws.Range("I129").Value = 2268.5715
$ws.Range("J129").Value = 2231.5
$ws.Range("K129").Value = 6805.7145
$ws.Range("L129").Value = 6694.5
$ws.Range("M129").Value = -1805.7145
$ws.Range("N129").Value = -16694.5
$ws.Range("H131").Value = 12822821
$ws.Range("J131").Value = 13891245
$ws.Range("L131").Value = 41673735
$ws.Range("N131").Value = -41683815

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 297.1905
$ws.Range("I2").Value = 219.5
$ws.Range("J2").Value = 367.81818
$ws.Range("K2").Value = 219.5
$ws.Range("L2").Value = 367.81818
$ws.Range("M2").Value = -106.5
$ws.Range("N2").Value = -593.81818
$ws.Range("H52").Value = 33333.332
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518
$ws.Range("H62").Value = 49999.668
$ws.Range("J62").Value = 49999
$ws.Range("L62").Value = 49999
$ws.Range("N62").Value = -51371
$ws.Range("H65").Value = 49999.668
$ws.Range("J65").Value = 49999
$ws.Range("L65").Value = 149997
$ws.Range("N65").Value = -156861
$ws.Range("H92").Value = 18749.8
$ws.Range("J92").Value = 18749.8
$ws.Range("L92").Value = 18749.8
$ws.Range("N92").Value = -22493.8
$ws.Range("H132").Value = 75850.55499999999
$ws.Range("I132").Value = 106625.52
$ws.Range("K132").Value = 319876.56
$ws.Range("M132").Value = -317346.56

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6250
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H22").Value = 987
$ws.Range("I22").Value = 861
$ws.Range("K22").Value = 861
$ws.Range("M22").Value = -566
$ws.Range("H27").Value = 987
$ws.Range("I27").Value = 861
$ws.Range("K27").Value = 861
$ws.Range("M27").Value = -754
$ws.Range("H55").Value = 581.8889
$ws.Range("J55").Value = 920.6
$ws.Range("L55").Value = 920.6
$ws.Range("N55").Value = -1266.6
$ws.Range("H132").Value = 2701.04
$ws.Range("I132").Value = 2661.647
$ws.Range("J132").Value = 2784.75
$ws.Range("K132").Value = 7984.941
$ws.Range("L132").Value = 8354.25
$ws.Range("M132").Value = -5454.941
$ws.Range("N132").Value = -13414.25
$ws.Range("H136").Value = 4322.909
$ws.Range("I136").Value = 2601.5144
$ws.Range("K136").Value = 7804.5432
$ws.Range("M136").Value = -5254.5432
$ws.Range("H139").Value = 88715
$ws.Range("J139").Value = 88715
$ws.Range("L139").Value = 88715
$ws.Range("N139").Value = -98995
$ws.Range("H141").Value = 112857.5
$ws.Range("J141").Value = 112857.5
$ws.Range("L141").Value = 112857.5
$ws.Range("N141").Value = -123217.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 4608.3335
$ws.Range("I100").Value = 609.05554
$ws.Range("K100").Value = 1218.11108
$ws.Range("M100").Value = -677.1110799999999
$ws.Range("H122").Value = 4143.278
$ws.Range("I122").Value = 1587.2
$ws.Range("J122").Value = 5126.385
$ws.Range("K122").Value = 4761.6
$ws.Range("L122").Value = 15379.155
$ws.Range("M122").Value = -2311.6
$ws.Range("N122").Value = -20279.155
$ws.Range("H132").Value = 4415.7646
$ws.Range("I132").Value = 1910.9
$ws.Range("J132").Value = 7994.143
$ws.Range("K132").Value = 5732.700000000001
$ws.Range("L132").Value = 23982.429
$ws.Range("M132").Value = -3202.700000000001
$ws.Range("N132").Value = -29042.429
$ws.Range("H136").Value = 4859.0967
$ws.Range("I136").Value = 3358.7778
$ws.Range("J136").Value = 6936.4614
$ws.Range("K136").Value = 10076.3334
$ws.Range("L136").Value = 20809.3842
$ws.Range("M136").Value = -7526.3334
$ws.Range("N136").Value = -25909.3842
$ws.Range("H138").Value = 97914.5
$ws.Range("L138").Value = 97429
$ws.Range("N138").Value = -107709
